$wb = $excel.ActiveWorkbook

function Update-HuanfeiSheet($ws) {
    # --- bump "want to go" counts on a few existing rows (F column) ---
    $ws.Range("F6").Value = 7149
    $ws.Range("F7").Value = 206
    $ws.Range("F8").Value = 172
    $ws.Range("F10").Value = 437

    # --- insert a brand-new event as the new row 11, pushing the old
    #     rows 11.. down by one ---
    $ws.Rows.Item(11).Insert()

    # Column A: numbered index, styled like the other index cells (copy
    # format from the row above, then set the literal number).
    $ws.Range("A10").Copy()
    $ws.Range("A11").PasteSpecial(-4122)
    $ws.Range("A11").Value = 10

    # Column B: a date-shaped string ("2024-06-30"). Force text storage
    # (so Excel doesn't silently turn it into a date serial), then copy
    # the plain formatting from a neighboring text cell so no stray
    # number format sticks around on the cell.
    $ws.Range("B11").NumberFormat = "@"
    $ws.Range("B11").Value = "2024-06-30"
    $ws.Range("C10").Copy()
    $ws.Range("B11").PasteSpecial(-4122)

    $ws.Range("C11").Value = "安徽·THO4·隙间皖韵之梦"
    $ws.Range("D11").Value = "北二环与新蚌埠路交汇处 蓝金湾大酒店"
    $ws.Range("E11").Value = "2024.06.30 10:00-06.30 17:00"
    $ws.Range("F11").Value = 0
    $ws.Range("G11").Value = 65
    $ws.Range("H11").Value = "https://show.bilibili.com/platform/detail.html?id=85119"
    $ws.Range("I11").Value = "//i2.hdslb.com/bfs/openplatform/202405/kuuarwvJ1714932457216.jpeg"

    # --- the rows that used to be 11/12/13 are now 12/13/14; their F
    #     (want-to-go) counts also changed in this refresh ---
    $ws.Range("F12").Value = 149
    $ws.Range("F13").Value = 191
    $ws.Range("F14").Value = 622
}

# Sheet 1 ("展览"): grows from A1:I13 to A1:I14.
$ws1 = $wb.Worksheets.Item(1)
Update-HuanfeiSheet($ws1)

# Sheet 4 ("全部类型"): same edits, but it already had one extra trailing
# row (the concert), so it grows from A1:I14 to A1:I15. That trailing row
# shifts automatically with the Rows.Insert() above and needs no further
# changes.
$ws4 = $wb.Worksheets.Item(4)
Update-HuanfeiSheet($ws4)
